{"js": "// Remove the last two rows (\"K\u1ebft lu\u1eadn\" and \"H\u01b0\u1edbng ph\u00e1t tri\u1ec3n\") from the\n// \"C\u00f4ng vi\u1ec7c \u0111\u00e3 l\u00e0m\" table, keeping row 13 (\"T\u1ea1o ch\u1ee9c n\u0103ng xu\u1ea5t phi\u1ebfu\n// nh\u1eadp xu\u1ea5t kho\") as the new last row.\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\n// Find the work-log table: the one whose first row text contains \"C\u00f4ng vi\u1ec7c \u0111\u00e3 l\u00e0m\"\nlet target = null;\nfor (let i = 0; i < tables.items.length; i++) {\n  const t = tables.items[i];\n  const rows = t.rows;\n  rows.load(\"items\");\n  await context.sync();\n  if (rows.items.length > 0) {\n    const firstRow = rows.items[0];\n    firstRow.load(\"values\");\n    await context.sync();\n    const headerText = (firstRow.values || []).flat().join(\" | \");\n    if (headerText.indexOf(\"C\u00f4ng vi\u1ec7c \u0111\u00e3 l\u00e0m\") !== -1) {\n      target = t;\n      break;\n    }\n  }\n}\n\nif (!target) {\n  throw new Error(\"Could not locate the 'C\u00f4ng vi\u1ec7c \u0111\u00e3 l\u00e0m' table\");\n}\n\nconst rows = target.rows;\nrows.load(\"items\");\nawait context.sync();\n\n// Identify the rows to remove by their text content, then delete them\n// from the highest index down so earlier indices stay valid.\nconst toDelete = [];\nfor (let i = 0; i < rows.items.length; i++) {\n  const row = rows.items[i];\n  row.load(\"values\");\n  await context.sync();\n  const text = (row.values || []).flat().join(\" | \").trim();\n  if (text.endsWith(\"K\u1ebft lu\u1eadn\") || text.endsWith(\"H\u01b0\u1edbng ph\u00e1t tri\u1ec3n\")) {\n    toDelete.push(i);\n  }\n}\n\ntoDelete.sort((a, b) => b - a);\nfor (const idx of toDelete) {\n  rows.items[idx].delete();\n}\nawait context.sync();\n", "ps1": "# Remove the last two rows (\"K\u1ebft lu\u1eadn\" and \"H\u01b0\u1edbng ph\u00e1t tri\u1ec3n\") from the\n# \"C\u00f4ng vi\u1ec7c \u0111\u00e3 l\u00e0m\" table, keeping row 13 (\"T\u1ea1o ch\u1ee9c n\u0103ng xu\u1ea5t phi\u1ebfu\n# nh\u1eadp xu\u1ea5t kho\") as the new last row.\n\n$d = $word.ActiveDocument\n\n# Locate the work-log table (header cell 2 reads \"C\u00f4ng vi\u1ec7c \u0111\u00e3 l\u00e0m\").\n$target = $null\nfor ($i = 1; $i -le $d.Tables.Count; $i++) {\n    $tbl = $d.Tables($i)\n    $headerText = $tbl.Cell(1, 2).Range.Text\n    if ($headerText -like \"*C\u00f4ng vi\u1ec7c \u0111\u00e3 l\u00e0m*\") {\n        $target = $tbl\n        break\n    }\n}\n\nif ($target -eq $null) {\n    throw \"Could not locate the 'C\u00f4ng vi\u1ec7c \u0111\u00e3 l\u00e0m' table\"\n}\n\n# Delete rows from the bottom up so indices of rows still to be removed\n# don't shift under us.\nfor ($i = $target.Rows.Count; $i -ge 1; $i--) {\n    $text = ($target.Cell($i, 2).Range.Text -replace \"[\\r\\a]\", \"\").Trim()\n    if ($text -eq \"K\u1ebft lu\u1eadn\" -or $text -eq \"H\u01b0\u1edbng ph\u00e1t tri\u1ec3n\") {\n        $target.Rows($i).Delete()\n    }\n}\n"}
